# "Settings back to qfin22 presentation - recalibration"
# Updates the recalibrated parameter values (mu, B, sig2, omega, alpha, beta,
# gamma, p, and the abs_epsi_autocorr series) across the AR / SETAR / GARCH /
# TARCH / AR_TARCH calibration sheets.
#
# NB: numeric literals are written in plain decimal (not scientific/E
# notation) because the host script parser does not accept an `E` exponent
# suffix on numeric literals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# AR sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = 0.000005686894670066777
$ws.Range("B3").Value = 0.7594194765332782
$ws.Range("B4").Value = 0.00003818422250271729
$ws.Range("B5").Value = "[1.0, 0.5024401975617963, 0.5170340989023596, 0.5151526414208251, 0.4593930042651452, 0.6282839412763306, 0.5223301499739756, 0.4766134253437698, 0.46935411377177305, 0.4847916378281686, 0.40884120712806094, 0.47009194025856643, 0.427483943005619, 0.44355259370347905, 0.4777091227681043, 0.39273836533563056, 0.373963972426112, 0.42831358332169517, 0.4445563893121485, 0.45163923791249444]"

# ---------------------------------------------------------------------------
# SETAR sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B4").Value = -0.0007412465094721902
$ws.Range("B5").Value = 0.7668115885867071
$ws.Range("B6").Value = 0.0000312521277313964
$ws.Range("B7").Value = 0.0004763676801595433
$ws.Range("B8").Value = 0.8427360108741354
$ws.Range("B9").Value = 0.0000246613095345031
$ws.Range("B10").Value = "[1.0, 0.48236897606498436, 0.4024134180534379, 0.4048342022802889, 0.41215414613972895, 0.5054159236088739, 0.44217273696953124, 0.4273123233828727, 0.49764603127382523, 0.4341678197555879, 0.4095480284918155, 0.4901852375775851, 0.3951967574584214, 0.35569711673820614, 0.38687532603101443, 0.35749159861734164, 0.4162739585216648, 0.34611572133412427, 0.4211265612902203, 0.3973175666397362]"

# ---------------------------------------------------------------------------
# GARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = -0.000002573746312684331
$ws.Range("B3").Value = 0.0000008680794102894604
$ws.Range("B4").Value = 0.2
$ws.Range("B5").Value = 0.78
$ws.Range("B6").Value = "[0.9999999999999998, 0.1265943927782083, 0.023792642029058945, 0.04259964721969289, 0.07569189153113443, 0.37570867857705087, 0.12205616030061298, 0.054589329705935066, 0.08065948718286206, 0.09004567668030015, 0.0528342912131126, 0.05218641023154138, 0.06931910234048162, 0.08084192064531064, 0.1600906921911198, 0.10137579125150088, 0.04754470287935678, 0.07330105497881605, 0.05464090813445139, 0.06857982238137578]"

# ---------------------------------------------------------------------------
# TARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.000002573746312684331
$ws.Range("B3").Value = 0.0000008680794102894604
$ws.Range("B4").Value = 0.2
$ws.Range("B5").Value = 0.68
$ws.Range("B6").Value = "[0.9999999999999999, 0.08811110765567585, -0.002881499357443855, 0.020988237332275374, 0.049581270212552535, 0.3564795851122364, 0.10146576105386855, 0.03660105635103391, 0.06641806723341273, 0.08389423926243869, 0.04304624257398931, 0.04341426898116409, 0.06246581010350186, 0.07694097961455972, 0.147186602900994, 0.09347098787671576, 0.036065812829057166, 0.07421764691341037, 0.05649178651135623, 0.05611763385758193]"
$ws.Range("B7").Value = 0.2

# ---------------------------------------------------------------------------
# AR_TARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = 0.000005686894670066851
$ws.Range("B3").Value = 0.0000007636844500543465
$ws.Range("B4").Value = 0.2
$ws.Range("B5").Value = 0.68
$ws.Range("B6").Value = "[1.0, 0.11392056853894429, 0.044512589345301494, 0.06282645517480463, 0.09015297350587634, 0.32657082054555725, 0.11148880772348989, 0.04292631727569616, 0.06603692206361411, 0.08998100805744644, 0.0676138663857636, 0.07465092306712333, 0.061681143721184446, 0.08707577237624899, 0.14585518375072298, 0.08350987900198734, 0.03947666657233857, 0.07469206187457274, 0.08866106585051738, 0.058740383911972595]"
$ws.Range("B7").Value = 0.2
$ws.Range("B9").Value = 0.759419476533278
